$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell so that it is stored as a *text* shared
# string even when the text looks like a number (e.g. "-666", "006666").
# Plain ".Value = ..." on such a string gets auto-converted to a numeric
# cell by the engine (just like real Excel would do when a user types a
# number into a cell), so for values that must remain text we either
# (a) copy an existing text cell that already holds the exact same string
#     (Copy + PasteSpecial values preserves the source cell's text type), or
# (b) stage the brand new string in a scratch cell far outside the used
#     range, force it to Text format, reset its style back to Normal so no
#     extra style is left referenced on the final cell, then copy that
#     scratch cell's value (still typed as text) into the destination and
#     clear the scratch cell.

function Set-TextFromExisting($destAddr, $sourceAddr) {
    $ws.Range($sourceAddr).Copy() | Out-Null
    $ws.Range($destAddr).PasteSpecial(-4163) | Out-Null
}

function Set-NewTextValue($destAddr, $text) {
    $scratch = $ws.Range("ZZ9000")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Style = "Normal"
    $scratch.Copy() | Out-Null
    $ws.Range($destAddr).PasteSpecial(-4163) | Out-Null
    $scratch.ClearContents() | Out-Null
}

# ---- Row 8 ----
$ws.Range("A8").Value = "Клемент Ворошилов"
Set-TextFromExisting "B8" "B4"                 # "-666"
Set-NewTextValue "C8" "02051945"
Set-TextFromExisting "D8" "D5"                 # "1945"
$ws.Range("E8").Value = "Винтовка Мосина"
$ws.Range("F8").Value = "1шт на взвод"
Set-TextFromExisting "G8" "G5"                 # "006666"
$ws.Range("H8").Value = "Не стреляет, но колит"
Set-TextFromExisting "I8" "I2"                 # "-"

# ---- Row 9 ----
$ws.Range("A9").Value = "Жуков Георгий"
Set-TextFromExisting "B9" "B4"                 # "-666"
Set-TextFromExisting "C9" "C5"                 # "09-05-1945"
Set-NewTextValue "D9" "2018"
$ws.Range("E9").Value = "Пленные Немцы"
$ws.Range("F9").Value = "40 полков"
Set-NewTextValue "G9" "001945"
$ws.Range("H9").Value = "Плохо работают"
Set-TextFromExisting "I9" "I2"                 # "-"
